$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "lista" column (B) for rows 86 to 121 from "Apruebo Dignidad" to "Unidad Constituyente"
$ws.Range("B86:B121").Value = "Unidad Constituyente"

# Update the selection in the sheet view to match the edited range
$ws.Range("B86:B121").Select()
